$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.441.51"
$ws.Range("E2").Value = "  -2.73%  "
$ws.Range("D3").Value = "1.803.52"
$ws.Range("E3").Value = "  -2.46%  "
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  +0.84%  "
$ws.Range("D5").Value = "1.007"
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("D6").Value = "307.93"
$ws.Range("E6").Value = "  -1.55%  "
$ws.Range("E7").Value = "  -1.41%  "
$ws.Range("D8").Value = "0.3654"
$ws.Range("E8").Value = "  -1.42%  "
$ws.Range("D9").Value = "0.07111"
$ws.Range("E9").Value = "  -2.23%  "
$ws.Range("D10").Value = "0.8747"
$ws.Range("D11").Value = "0.07730"
$ws.Range("E11").Value = "  -1.17%  "
$ws.Range("D12").Value = "19.31"
$ws.Range("E12").Value = "  -3.34%  "
$ws.Range("D13").Value = "1.815.62"
$ws.Range("E13").Value = "  -3.39%  "
$ws.Range("D14").Value = "5.263"
$ws.Range("E14").Value = "  -1.92%  "
$ws.Range("D15").Value = "6.338"
$ws.Range("E15").Value = "  -2.43%  "
$ws.Range("D16").Value = "85.83"
$ws.Range("E16").Value = "  -5.92%  "
$ws.Range("E17").Value = "  +0.84%  "
$ws.Range("D18").Value = "0.000008569"
$ws.Range("E18").Value = "  -3.95%  "
$ws.Range("E19").Value = "  +0.58%  "
$ws.Range("D20").Value = "26.498.41"
$ws.Range("E20").Value = "  -2.63%  "
$ws.Range("D21").Value = "14.21"
$ws.Range("E21").Value = "  -3.22%  "
$ws.Range("D22").Value = "4.966"
$ws.Range("E22").Value = "  -1.74%  "
$ws.Range("E23").Value = "  -0.81%  "
$ws.Range("D24").Value = "1.983"
$ws.Range("E24").Value = "  +1.34%  "
$ws.Range("D25").Value = "150.40"
$ws.Range("E25").Value = "  -0.89%  "
$ws.Range("D26").Value = "17.90"
$ws.Range("E26").Value = "  -2.57%  "
$ws.Range("D27").Value = "1.999"
$ws.Range("E27").Value = "  -3.01%  "
$ws.Range("D28").Value = "112.45"
$ws.Range("E28").Value = "  -2.66%  "
$ws.Range("D29").Value = "4.810"
$ws.Range("E29").Value = "  -4.77%  "
$ws.Range("D30").Value = "0.08646"
$ws.Range("E30").Value = "  -2.01%  "
$ws.Range("D31").Value = "3.049"
$ws.Range("E31").Value = "  -1.50%  "
$ws.Range("D32").Value = "0.7269"
$ws.Range("E32").Value = "  -4.47%  "
$ws.Range("D33").Value = "4.423"
$ws.Range("E33").Value = "  -1.50%  "
$ws.Range("E34").Value = "  -4.76%  "
$ws.Range("D35").Value = "1.009"
$ws.Range("E35").Value = "  +0.81%  "
$ws.Range("D36").Value = "2.555"
$ws.Range("E36").Value = "  -6.45%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "0.01926"
$ws.Range("E38").Value = "  -0.63%  "
$ws.Range("D39").Value = "0.05083"
$ws.Range("E39").Value = "  -2.89%  "
$ws.Range("D40").Value = "2.879"
$ws.Range("E40").Value = "  -1.98%  "
$ws.Range("D41").Value = "6.930"
$ws.Range("E41").Value = "  -1.92%  "
$ws.Range("D42").Value = "0.4993"
$ws.Range("E42").Value = "  -1.86%  "
$ws.Range("D43").Value = "0.1561"
$ws.Range("E43").Value = "  -3.97%  "
$ws.Range("D44").Value = "8.086"
$ws.Range("E44").Value = "  -3.28%  "
$ws.Range("D45").Value = "1.008"
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("D46").Value = "0.4590"
$ws.Range("E46").Value = "  -3.98%  "
$ws.Range("D47").Value = "9.952"
$ws.Range("E47").Value = "  -3.82%  "
$ws.Range("D48").Value = "101.72"
$ws.Range("D49").Value = "1.586"
$ws.Range("E49").Value = "  -2.81%  "
$ws.Range("D50").Value = "0.05987"
$ws.Range("E50").Value = "  -3.67%  "
$ws.Range("D51").Value = "63.73"
$ws.Range("E51").Value = "  -2.80%  "
